# Auto-generated edit script applying Pandaemonium_Profits.xlsx diff
# Updates currentAveragePrice/currentAveragePriceNQ/currentAveragePriceHQ/
# LevePriceNQ/LevePriceHQ/LeveProfitNQ/LeveProfitHQ (columns H-N) for
# various leve rows across all eight crafting-class sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1512.909
$ws.Cells.Item(38, 10).Value = 1929.375
$ws.Cells.Item(38, 12).Value = 5788.125
$ws.Cells.Item(38, 14).Value = -6532.125
$ws.Cells.Item(109, 8).Value = 60800
$ws.Cells.Item(109, 10).Value = 60800
$ws.Cells.Item(109, 12).Value = 60800
$ws.Cells.Item(109, 14).Value = -63574
$ws.Cells.Item(129, 8).Value = 973.2037
$ws.Cells.Item(129, 9).Value = 429.16666
$ws.Cells.Item(129, 10).Value = 1041.2084
$ws.Cells.Item(129, 11).Value = 1287.49998
$ws.Cells.Item(129, 12).Value = 3123.6252
$ws.Cells.Item(129, 13).Value = 3712.50002
$ws.Cells.Item(129, 14).Value = -13123.6252
$ws.Cells.Item(137, 8).Value = 928717.75
$ws.Cells.Item(137, 9).Value = 3000.4119
$ws.Cells.Item(137, 10).Value = 1756991.2
$ws.Cells.Item(137, 11).Value = 9001.235700000001
$ws.Cells.Item(137, 12).Value = 5270973.6
$ws.Cells.Item(137, 13).Value = -6451.235700000001
$ws.Cells.Item(137, 14).Value = -5276073.6

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1625.6364
$ws.Cells.Item(45, 9).Value = 1572.3125
$ws.Cells.Item(45, 11).Value = 1572.3125
$ws.Cells.Item(45, 13).Value = -1195.3125
$ws.Cells.Item(61, 8).Value = 5169.2983
$ws.Cells.Item(61, 9).Value = 3365.1633
$ws.Cells.Item(61, 10).Value = 10080.556
$ws.Cells.Item(61, 11).Value = 3365.1633
$ws.Cells.Item(61, 12).Value = 10080.556
$ws.Cells.Item(61, 13).Value = -3153.1633
$ws.Cells.Item(61, 14).Value = -10504.556
$ws.Cells.Item(63, 8).Value = 3772
$ws.Cells.Item(63, 9).Value = 2901
$ws.Cells.Item(63, 10).Value = 4933.3335
$ws.Cells.Item(63, 11).Value = 2901
$ws.Cells.Item(63, 12).Value = 4933.3335
$ws.Cells.Item(63, 13).Value = -2215
$ws.Cells.Item(63, 14).Value = -6305.3335
$ws.Cells.Item(66, 8).Value = 3772
$ws.Cells.Item(66, 9).Value = 2901
$ws.Cells.Item(66, 10).Value = 4933.3335
$ws.Cells.Item(66, 11).Value = 14505
$ws.Cells.Item(66, 12).Value = 24666.6675
$ws.Cells.Item(66, 13).Value = -11073
$ws.Cells.Item(66, 14).Value = -31530.6675
$ws.Cells.Item(74, 8).Value = 4859.475
$ws.Cells.Item(74, 9).Value = 2048.2856
$ws.Cells.Item(74, 11).Value = 2048.2856
$ws.Cells.Item(74, 13).Value = -1174.2856
$ws.Cells.Item(77, 8).Value = 4859.475
$ws.Cells.Item(77, 9).Value = 2048.2856
$ws.Cells.Item(77, 11).Value = 10241.428
$ws.Cells.Item(77, 13).Value = -5873.428
$ws.Cells.Item(94, 8).Value = 39000
$ws.Cells.Item(94, 10).Value = 39000
$ws.Cells.Item(94, 12).Value = 39000
$ws.Cells.Item(94, 14).Value = -40802
$ws.Cells.Item(97, 8).Value = 1531.4286
$ws.Cells.Item(97, 9).Value = 1620
$ws.Cells.Item(97, 10).Value = 1000
$ws.Cells.Item(97, 11).Value = 1620
$ws.Cells.Item(97, 12).Value = 1000
$ws.Cells.Item(97, 13).Value = -1124
$ws.Cells.Item(97, 14).Value = -1992
$ws.Cells.Item(136, 8).Value = 5169.2983
$ws.Cells.Item(136, 9).Value = 3365.1633
$ws.Cells.Item(136, 10).Value = 10080.556
$ws.Cells.Item(136, 11).Value = 10095.4899
$ws.Cells.Item(136, 12).Value = 30241.668
$ws.Cells.Item(136, 13).Value = -7545.4899
$ws.Cells.Item(136, 14).Value = -35341.66800000001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(103, 8).Value = 30657
$ws.Cells.Item(103, 10).Value = 30657
$ws.Cells.Item(103, 12).Value = 30657
$ws.Cells.Item(103, 14).Value = -33001
$ws.Cells.Item(105, 8).Value = 5248.125
$ws.Cells.Item(105, 9).Value = 5536.154
$ws.Cells.Item(105, 11).Value = 5536.154
$ws.Cells.Item(105, 13).Value = -3789.154
$ws.Cells.Item(106, 8).Value = 18671
$ws.Cells.Item(106, 10).Value = 18671
$ws.Cells.Item(106, 12).Value = 18671
$ws.Cells.Item(106, 14).Value = -21195

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 630936.5
$ws.Cells.Item(31, 9).Value = 4910.76
$ws.Cells.Item(31, 10).Value = 1091249.5
$ws.Cells.Item(31, 11).Value = 4910.76
$ws.Cells.Item(31, 12).Value = 1091249.5
$ws.Cells.Item(31, 13).Value = -4615.76
$ws.Cells.Item(31, 14).Value = -1091839.5
$ws.Cells.Item(34, 8).Value = 630936.5
$ws.Cells.Item(34, 9).Value = 4910.76
$ws.Cells.Item(34, 10).Value = 1091249.5
$ws.Cells.Item(34, 11).Value = 4910.76
$ws.Cells.Item(34, 12).Value = 1091249.5
$ws.Cells.Item(34, 13).Value = -4708.76
$ws.Cells.Item(34, 14).Value = -1091653.5
$ws.Cells.Item(35, 8).Value = 4225
$ws.Cells.Item(35, 9).Value = 450
$ws.Cells.Item(35, 11).Value = 450
$ws.Cells.Item(35, 13).Value = -156
$ws.Cells.Item(62, 8).Value = 3420.5557
$ws.Cells.Item(62, 9).Value = 3464.1667
$ws.Cells.Item(62, 10).Value = 3333.3333
$ws.Cells.Item(62, 11).Value = 3464.1667
$ws.Cells.Item(62, 12).Value = 3333.3333
$ws.Cells.Item(62, 13).Value = -2840.1667
$ws.Cells.Item(62, 14).Value = -4581.3333
$ws.Cells.Item(65, 8).Value = 3420.5557
$ws.Cells.Item(65, 9).Value = 3464.1667
$ws.Cells.Item(65, 10).Value = 3333.3333
$ws.Cells.Item(65, 11).Value = 17320.8335
$ws.Cells.Item(65, 12).Value = 16666.6665
$ws.Cells.Item(65, 13).Value = -14200.8335
$ws.Cells.Item(65, 14).Value = -22906.6665
$ws.Cells.Item(95, 8).Value = 22449.6
$ws.Cells.Item(95, 10).Value = 22449.6
$ws.Cells.Item(95, 12).Value = 22449.6
$ws.Cells.Item(95, 14).Value = -27941.6
$ws.Cells.Item(105, 8).Value = 801.4737
$ws.Cells.Item(105, 9).Value = 823.2222
$ws.Cells.Item(105, 10).Value = 410
$ws.Cells.Item(105, 11).Value = 823.2222
$ws.Cells.Item(105, 12).Value = 410
$ws.Cells.Item(105, 13).Value = 923.7778
$ws.Cells.Item(105, 14).Value = -3904
$ws.Cells.Item(107, 8).Value = 1212.0834
$ws.Cells.Item(107, 9).Value = 1276.8182
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 11).Value = 1276.8182
$ws.Cells.Item(107, 12).Value = 500
$ws.Cells.Item(107, 13).Value = 643.1818000000001
$ws.Cells.Item(107, 14).Value = -4340

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 91.5
$ws.Cells.Item(33, 9).Value = 71.40000000000001
$ws.Cells.Item(33, 11).Value = 428.4
$ws.Cells.Item(33, 13).Value = -145.4
$ws.Cells.Item(107, 8).Value = 723.3538
$ws.Cells.Item(107, 9).Value = 276.33334
$ws.Cells.Item(107, 10).Value = 1729.15
$ws.Cells.Item(107, 11).Value = 829.0000200000001
$ws.Cells.Item(107, 12).Value = 5187.450000000001
$ws.Cells.Item(107, 13).Value = 1090.99998
$ws.Cells.Item(107, 14).Value = -9027.450000000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(130, 8).Value = 69266.664
$ws.Cells.Item(130, 10).Value = 69266.664
$ws.Cells.Item(130, 12).Value = 69266.664
$ws.Cells.Item(130, 14).Value = -79306.664

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 40406.285
$ws.Cells.Item(133, 10).Value = 40406.285
$ws.Cells.Item(133, 12).Value = 40406.285
$ws.Cells.Item(133, 14).Value = -45466.285
$ws.Cells.Item(136, 8).Value = 3002
$ws.Cells.Item(136, 9).Value = 2741.087
$ws.Cells.Item(136, 10).Value = 6002.5
$ws.Cells.Item(136, 11).Value = 8223.261
$ws.Cells.Item(136, 12).Value = 18007.5
$ws.Cells.Item(136, 13).Value = -5673.261
$ws.Cells.Item(136, 14).Value = -23107.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2000
$ws.Cells.Item(96, 9).Value = 2000
$ws.Cells.Item(96, 11).Value = 2000
$ws.Cells.Item(96, 13).Value = -627
$ws.Cells.Item(112, 8).Value = 79800
$ws.Cells.Item(112, 10).Value = 79800
$ws.Cells.Item(112, 12).Value = 79800
$ws.Cells.Item(112, 14).Value = -82754
$ws.Cells.Item(136, 8).Value = 4968.553
$ws.Cells.Item(136, 9).Value = 1511.9584
$ws.Cells.Item(136, 10).Value = 8575.434999999999
$ws.Cells.Item(136, 11).Value = 4535.8752
$ws.Cells.Item(136, 12).Value = 25726.305
$ws.Cells.Item(136, 13).Value = -1985.8752
$ws.Cells.Item(136, 14).Value = -30826.305
